$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2080.8667
$ws.Range("I106").Value = 1763.091
$ws.Range("K106").Value = 1763.091
$ws.Range("M106").Value = -1132.091
$ws.Range("H112").Value = 1649
$ws.Range("J112").Value = 1681.5652
$ws.Range("L112").Value = 5044.6956
$ws.Range("N112").Value = -7260.6956
$ws.Range("H125").Value = 2941806
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2941806
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 26476254
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -26481174
$ws.Range("H137").Value = 1219.4667
$ws.Range("I137").Value = 1285.5385
$ws.Range("J137").Value = 790
$ws.Range("K137").Value = 3856.6155
$ws.Range("L137").Value = 2370
$ws.Range("M137").Value = -1306.6155
$ws.Range("N137").Value = -7470
$ws.Range("H138").Value = 3594.701
$ws.Range("I138").Value = 1619.8667
$ws.Range("J138").Value = 5303.6924
$ws.Range("K138").Value = 4859.6001
$ws.Range("L138").Value = 15911.0772
$ws.Range("M138").Value = 280.3999000000003
$ws.Range("N138").Value = -26191.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7214.2656
$ws.Range("I32").Value = 6440.6104
$ws.Range("J32").Value = 37000
$ws.Range("K32").Value = 6440.6104
$ws.Range("L32").Value = 37000
$ws.Range("M32").Value = -6153.6104
$ws.Range("N32").Value = -37574
$ws.Range("H74").Value = 1424.8889
$ws.Range("I74").Value = 1017.625
$ws.Range("J74").Value = 2239.4167
$ws.Range("K74").Value = 1017.625
$ws.Range("L74").Value = 2239.4167
$ws.Range("M74").Value = -143.625
$ws.Range("N74").Value = -3987.4167
$ws.Range("H77").Value = 1424.8889
$ws.Range("I77").Value = 1017.625
$ws.Range("J77").Value = 2239.4167
$ws.Range("K77").Value = 5088.125
$ws.Range("L77").Value = 11197.0835
$ws.Range("M77").Value = -720.125
$ws.Range("N77").Value = -19933.0835
$ws.Range("H97").Value = 607.37933
$ws.Range("I97").Value = 596.26086
$ws.Range("J97").Value = 650
$ws.Range("K97").Value = 596.26086
$ws.Range("L97").Value = 650
$ws.Range("M97").Value = -100.26086
$ws.Range("N97").Value = -1642
$ws.Range("H133").Value = 82298.664
$ws.Range("J133").Value = 82298.664
$ws.Range("L133").Value = 82298.664
$ws.Range("N133").Value = -87358.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 527.5833
$ws.Range("I94").Value = 527.5833
$ws.Range("K94").Value = 527.5833
$ws.Range("M94").Value = -76.58330000000001
$ws.Range("H99").Value = 1680.909
$ws.Range("I99").Value = 1125
$ws.Range("K99").Value = 1125
$ws.Range("M99").Value = 373
$ws.Range("H107").Value = 649
$ws.Range("I107").Value = 649
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 649
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1271
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3225.1304
$ws.Range("I31").Value = 1368.0952
$ws.Range("J31").Value = 7259.3794
$ws.Range("K31").Value = 1368.0952
$ws.Range("L31").Value = 7259.3794
$ws.Range("M31").Value = -1073.0952
$ws.Range("N31").Value = -7849.3794
$ws.Range("H34").Value = 3225.1304
$ws.Range("I34").Value = 1368.0952
$ws.Range("J34").Value = 7259.3794
$ws.Range("K34").Value = 1368.0952
$ws.Range("L34").Value = 7259.3794
$ws.Range("M34").Value = -1166.0952
$ws.Range("N34").Value = -7663.3794
$ws.Range("H122").Value = 1254308.2
$ws.Range("I122").Value = 3012
$ws.Range("J122").Value = 1671407
$ws.Range("K122").Value = 9036
$ws.Range("L122").Value = 5014221
$ws.Range("M122").Value = -6586
$ws.Range("N122").Value = -5019121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 134.47058
$ws.Range("I2").Value = 433.6
$ws.Range("J2").Value = 9.833333
$ws.Range("K2").Value = 2601.6
$ws.Range("L2").Value = 58.999998
$ws.Range("M2").Value = -2488.6
$ws.Range("N2").Value = -284.999998
$ws.Range("H105").Value = 8757.143
$ws.Range("J105").Value = 8757.143
$ws.Range("L105").Value = 26271.429
$ws.Range("N105").Value = -31513.429
$ws.Range("H107").Value = 592837.4399999999
$ws.Range("J107").Value = 756761.5600000001
$ws.Range("L107").Value = 2270284.68
$ws.Range("N107").Value = -2274124.68
$ws.Range("H113").Value = 1871.6666
$ws.Range("I113").Value = 2163.5715
$ws.Range("J113").Value = 850
$ws.Range("K113").Value = 6490.7145
$ws.Range("L113").Value = 2550
$ws.Range("M113").Value = -4320.7145
$ws.Range("N113").Value = -6890
$ws.Range("H131").Value = 857.1900000000001
$ws.Range("I131").Value = 567.25
$ws.Range("J131").Value = 869.2708
$ws.Range("K131").Value = 1701.75
$ws.Range("L131").Value = 2607.8124
$ws.Range("M131").Value = 3338.25
$ws.Range("N131").Value = -12687.8124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 974.8261
$ws.Range("I97").Value = 986.8182
$ws.Range("J97").Value = 711
$ws.Range("K97").Value = 986.8182
$ws.Range("L97").Value = 711
$ws.Range("M97").Value = -490.8182
$ws.Range("N97").Value = -1703
$ws.Range("H132").Value = 1596.0339
$ws.Range("I132").Value = 1318.5946
$ws.Range("J132").Value = 2062.6365
$ws.Range("K132").Value = 3955.7838
$ws.Range("L132").Value = 6187.9095
$ws.Range("M132").Value = -1425.7838
$ws.Range("N132").Value = -11247.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 59079.945
$ws.Range("I7").Value = 74602.86
$ws.Range("K7").Value = 74602.86
$ws.Range("M7").Value = -74490.86
$ws.Range("H45").Value = 32999
$ws.Range("J45").Value = 32999
$ws.Range("L45").Value = 32999
$ws.Range("N45").Value = -33813
$ws.Range("H122").Value = 6538075.5
$ws.Range("I122").Value = 12346977
$ws.Range("J122").Value = 3060.625
$ws.Range("K122").Value = 37040931
$ws.Range("L122").Value = 9181.875
$ws.Range("M122").Value = -37038481
$ws.Range("N122").Value = -14081.875
$ws.Range("H126").Value = 59079.945
$ws.Range("I126").Value = 74602.86
$ws.Range("K126").Value = 223808.58
$ws.Range("M126").Value = -221338.58
$ws.Range("H132").Value = 7071.164
$ws.Range("I132").Value = 7274.4
$ws.Range("J132").Value = 6499.5625
$ws.Range("K132").Value = 21823.2
$ws.Range("L132").Value = 19498.6875
$ws.Range("M132").Value = -19293.2
$ws.Range("N132").Value = -24558.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 578.86664
$ws.Range("I107").Value = 549.4167
$ws.Range("J107").Value = 696.6667
$ws.Range("K107").Value = 1648.2501
$ws.Range("L107").Value = 2090.0001
$ws.Range("M107").Value = 271.7499
$ws.Range("N107").Value = -5930.0001
$ws.Range("H122").Value = 33972.29
$ws.Range("I122").Value = 45032.434
$ws.Range("J122").Value = 2174.375
$ws.Range("K122").Value = 135097.302
$ws.Range("L122").Value = 6523.125
$ws.Range("M122").Value = -132647.302
$ws.Range("N122").Value = -11423.125
$ws.Range("H132").Value = 1348.4333
$ws.Range("I132").Value = 893.26666
$ws.Range("J132").Value = 2713.9333
$ws.Range("K132").Value = 2679.79998
$ws.Range("L132").Value = 8141.7999
$ws.Range("M132").Value = -149.7999799999998
$ws.Range("N132").Value = -13201.7999
